$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 501 (Fecha 45077), pushing every
# existing record from row 501 down one row (old row 501 -> new row 502, ...,
# old row 573 -> new row 574). Insert a blank row first so everything below
# shifts down intact.
$ws.Rows.Item(501).Insert()

# The freshly inserted row 501 is blank; seed it with the same "shape" as the
# row that just got pushed down to 502, then overwrite the four fields that
# actually differ for the new record (Fecha, Volumen, Precio promedio
# ponderado, Precio $/Kg).
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(501, $col).Value = $ws.Cells.Item(502, $col).Value()
}

$ws.Cells.Item(501, 4).Value  = 45077   # D501 - Fecha
$ws.Cells.Item(501, 10).Value = 175     # J501 - Volumen
$ws.Cells.Item(501, 13).Value = 4186    # M501 - Precio promedio ponderado
$ws.Cells.Item(501, 16).Value = 1395    # P501 - Precio $/Kg
